# "Se agrego boton para agregar filas manualmente"
# Append new invoice rows (6-13) below the existing data, mirroring the
# behaviour of a manual "add row" button: every new cell is entered as
# literal text (so dates/numbers aren't auto-converted by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("SKM_C335124092312180_0026.jpg", "12/07/2024", "2",     "200,00", "165,29", "21,00"),
    @("SKM_C335124092312180_0033.jpg", "18/07/2024", "11307", "170,01", "140,50", "21"),
    @("SKM_C335124092312180_0032.jpg", "28/08/2024", "13536", "150,00", "123,97", "21,00"),
    @("SKM_C335124092312180_0032.jpg", "28/08/2024", "3.536", "150,00", "123,97", "21,00"),
    @("SKM_C335124092312180_0030.jpg", "16/08/2024", "12980", "200,00", "165,29", "21,00"),
    @("SKM_C335124092312180_0031.jpg", "21/08/2024", "6938R", "150,00", "123,97", "21,00"),
    @("SKM_C335124092312180_0028.jpg", "29/07/2024", "9838R", "47,06",  "38,89",  "8,17"),
    @("",                              "01/01/2024", "9218",  "100,21", "12,2",   "21")
)

$startRow = 6
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]

    # Force every column to Text format first so values like "12/07/2024"
    # or "2" are stored as literal strings, not dates/numbers.
    $rowRange = $ws.Range("A" + $r + ":F" + $r)
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}
